$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update bought-item label (Vietnamese) to be upper-case
$ws.Range("C28").Value = "ĐÃ MUA: "

# Re-word the bet-coin prompts (English + Vietnamese)
$ws.Range("B31").Value = "INTER YOUR STAKE: "
$ws.Range("C31").Value = "NHẬP SỐ TIỀN CƯỢC: "

# Fix spelling of the win message (English)
$ws.Range("B29").Value = "CONGRATULATIONS! YOU WIN!!!"

# Move the active selection to reflect the latest edit position
$ws.Range("B29").Select()
